$d = $word.ActiveDocument

# ----------------------------------------------------------------------
# Edit 1: " ... performed on the data values " -> " ... performed on the
# data " + "point" + bookmark(_GoBack) + " " (+ existing "with the
# smallest length measurement" continues unchanged).
# ----------------------------------------------------------------------
$find1 = $d.Content.Find
$find1.Execute("performed on the data values ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($find1.Found) {
    $whole = $find1.Parent
    $cutStart = $whole.End - 7   # start of the trailing "values " (7 chars incl. space)

    # Drop the trailing "values " - a pure tail-delete so the remaining
    # run ("... on the data ") keeps its original run identity.
    $sub = $d.Range($cutStart, $whole.End)
    $sub.Delete()

    # Insert "point" as a new run right after, forcing a run split by
    # briefly diverging formatting, then cleanly resetting it.
    $p1 = $d.Range($cutStart, $cutStart)
    $p1.InsertAfter("point")
    $p1.Font.Bold = 1

    # Insert the trailing " " right after "point" (still diverging).
    $afterPoint = $cutStart + 5
    $p2 = $d.Range($afterPoint, $afterPoint)
    $p2.InsertAfter(" ")
    $p2.Font.Bold = 1

    # Re-home the document's sole "_GoBack" bookmark to sit exactly
    # between "point" and the following " " run (collapsed range).
    $bmRange = $d.Range($afterPoint, $afterPoint)
    $d.Bookmarks.Add("_GoBack", $bmRange)

    # Now clear the temporary bold divergence across the whole span in
    # one shot so "point" and " " resolve to clean (rsid-less) runs.
    $resetRange = $d.Range($cutStart, $afterPoint + 1)
    $resetRange.Font.Bold = 0
}

# ----------------------------------------------------------------------
# Edit 2: "Use the smallest measured wavelength to estimate largest
# possible error." -> "...wavelength " + "data point " + "to estimate...".
# ----------------------------------------------------------------------
$find2 = $d.Content.Find
$find2.Execute("Use the smallest measured wavelength to estimate largest possible error.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($find2.Found) {
    $whole2 = $find2.Parent
    $prefixLen2 = "Use the smallest measured wavelength ".Length
    $cutStart2 = $whole2.Start + $prefixLen2

    $sub2 = $d.Range($cutStart2, $whole2.End)
    $suffixText2 = $sub2.Text
    $sub2.Delete()

    # Re-insert the suffix, unmarked, directly after the (now shortened)
    # prefix run.
    $p2b = $d.Range($cutStart2, $cutStart2)
    $p2b.InsertAfter($suffixText2)

    # Collapse to the start of the just-reinserted suffix and insert
    # "data point " ahead of it - this reliably produces 3 distinct runs.
    $p2b.Collapse(1)
    $p2b.InsertBefore("data point ")
    $p2b.Font.Bold = 1
    $p2b.Font.Bold = 0
}

# ----------------------------------------------------------------------
# Edit 3: "Use the smallest measured distance and time to estimate
# largest possible error." -> "...time " + "data point " + "to estimate...".
# ----------------------------------------------------------------------
$find3 = $d.Content.Find
$find3.Execute("Use the smallest measured distance and time to estimate largest possible error.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($find3.Found) {
    $whole3 = $find3.Parent
    $prefixLen3 = "Use the smallest measured distance and time ".Length
    $cutStart3 = $whole3.Start + $prefixLen3

    $sub3 = $d.Range($cutStart3, $whole3.End)
    $suffixText3 = $sub3.Text
    $sub3.Delete()

    $p3b = $d.Range($cutStart3, $cutStart3)
    $p3b.InsertAfter($suffixText3)

    $p3b.Collapse(1)
    $p3b.InsertBefore("data point ")
    $p3b.Font.Bold = 1
    $p3b.Font.Bold = 0
}

# ----------------------------------------------------------------------
# Edit 4 (the old "_GoBack" bookmark near "... quiet room.") is handled
# implicitly: Word (and this host) only ever keeps a single "_GoBack"
# bookmark, so re-adding it in Edit 1 above already relocated it away
# from its old position.
# ----------------------------------------------------------------------

Write-Output "done"
